$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -------------
# The shared string "Ready for handoff" is referenced from every localized
# status cell on all three sheets. Every one of those cells must be
# rewritten to the new text so the old shared-string entry becomes
# unreferenced (and is dropped on save) while all updated cells collapse
# onto a single new "In Translation" shared string, matching the diff.

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Range("E4").Value = "In Translation"
$ws1.Range("F4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("C4").Value = "In Translation"

# --- Column width changes ----------------------------------------------
# Overview!E:F and the Status column (C) on the zh-cn/de-de sheets shrink
# (their text got shorter: "Ready for handoff" -> "In Translation").
# ColumnWidth is specified in characters; set the nearest width to the
# narrower target.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
